# Edit: slide 1, "TextBox 7" shape (GVHD/SV info box) -- enlarge font sizes,
# add tab characters, and resize/reposition the textbox.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(4)   # "TextBox 7"
$tf = $shp.TextFrame
$tr = $tf.TextRange

# --- Paragraph 1: "SV/ Nhóm SV" + ": 1. Hồ Tuấn Phước " ---
$para1 = $tr.Paragraphs(1,1)

$p1r1 = $para1.Runs(1,1)
$p1r1.Font.Size = 28

$p1r2 = $para1.Runs(2,1)
$p1r2.Font.Size = 28
$p1r2.Text = ": `t1. Hồ Tuấn Phước "

# --- Paragraph 2: tab + "  " + "2. " + "Thái Nguyễn Thiện Duyên" ---
$para2 = $tr.Paragraphs(2,1)

$p2r1 = $para2.Runs(1,1)
$p2r1.Font.Size = 28

$p2r2 = $para2.Runs(2,1)
$p2r2.Font.Size = 28
$p2r2.Text = "  `t`t"

$p2r3 = $para2.Runs(3,1)
$p2r3.Font.Size = 28

$p2r4 = $para2.Runs(4,1)
$p2r4.Font.Size = 24

# --- Resize / reposition the textbox (EMU -> pt, 12700 EMU per pt) ---
$shp.Left = 188.75
$shp.Top = 289.25
$shp.Width = 652.30002
$shp.Height = 75.05
